# Apply crypto price/volume/coin updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column D (Price) to be treated as text so values like "30.815.85"
# or "1.002" are stored verbatim instead of being parsed as numbers/dates.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.815.85'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '1.889.18'
$ws.Range("E3").Value = '  +2.68%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '237.78'
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.4794'
$ws.Range("E7").Value = '  +2.70%  '
$ws.Range("D8").Value = '0.2858'
$ws.Range("E8").Value = '  +5.51%  '
$ws.Range("D9").Value = '0.06525'
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").Value = '18.82'
$ws.Range("E10").Value = '  +16.98%  '
$ws.Range("D11").Value = '1.887.34'
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").Value = '0.07573'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '94.94'
$ws.Range("E13").Value = '  +13.48%  '
$ws.Range("D14").Value = '5.136'
$ws.Range("E14").Value = '  +4.21%  '
$ws.Range("D15").Value = '0.6529'
$ws.Range("E15").Value = '  +5.57%  '
$ws.Range("D16").Value = '297.40'
$ws.Range("E16").Value = '  +31.68%  '
$ws.Range("D17").Value = '30.828.04'
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("D18").Value = '13.19'
$ws.Range("E18").Value = '  +6.84%  '
$ws.Range("D19").Value = '0.9999'
$ws.Range("D20").Value = '0.000007523'
$ws.Range("E20").Value = '  +3.48%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.119.10'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.178'
$ws.Range("E23").Value = '  +5.97%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '6.148'
$ws.Range("E24").Value = '  +5.19%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.316'
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '168.10'
$ws.Range("E26").Value = '  +2.39%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '19.69'
$ws.Range("E27").Value = '  +11.05%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.953'
$ws.Range("E28").Value = '  +5.09%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.1066'
$ws.Range("E29").Value = '  +3.33%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.362'
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.196'
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.970'
$ws.Range("E32").Value = '  +4.25%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.05007'
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.174'
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7251'
$ws.Range("E35").Value = '  +2.58%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.717'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.01948'
$ws.Range("E37").Value = '  +4.38%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.723'
$ws.Range("E38").Value = '  +2.75%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.056'
$ws.Range("E39").Value = '  +7.20%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.8968'
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '107.63'
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4195'
$ws.Range("E43").Value = '  +4.71%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.596'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '66.05'
$ws.Range("E45").Value = '  +10.63%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.344'
$ws.Range("E46").Value = '  +4.63%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1226'
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.881'
$ws.Range("E48").Value = '  +3.21%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '34.62'
$ws.Range("E49").Value = '  +5.44%  '
$ws.Range("D50").Value = '0.05629'
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.394'
$ws.Range("E51").Value = '  +2.73%  '
